$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 76
$ws.Range("I2").Value = 189
$ws.Range("J2").Value = 773
$ws.Range("K2").Value = 6
$ws.Range("L2").Value = 205
$ws.Range("M2").Value = 19
$ws.Range("N2").Value = 131
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 13
$ws.Range("S2").Value = 104
$ws.Range("T2").Value = 163
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 1176
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1196
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 8
